$d = $word.ActiveDocument

# 1. "CSCI UA.0060 Fall 2024" -> "CSCI UA.0060 Spring 2025"
$d.Content.Find.Execute("Fall 2024", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Spring 2025", 2)

# 2. "...comfortable so long as..." -> "...comfortable with so long as..."
$d.Content.Find.Execute("comfortable so long as", $true, $false, $false, $false, $false,
                         $true, 1, $false, "comfortable with so long as", 2)
